$d = $word.ActiveDocument

# --- Locate the paragraph that holds the "{m:commentblock ...}" field token ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "{m:commentblock*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the {m:commentblock...} paragraph"
}

$pStart = $target.Range.Start

# Helper: split off the sub-range [startPos, endPos) into its own run by
# toggling Bold on then off again. The host only materializes a new run
# boundary once direct formatting has actually been touched, so the
# flip-flop leaves the visible formatting untouched while forcing the split.
function Split-Range($startPos, $endPos) {
    $r = $d.Range($startPos, $endPos)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# --- 1) Insert the missing space right after the ':' (before "comment") ---
# Original text:  {m:commentblock some important comment}
# Target text:    {m: commentblock some important comment}
$insertPos = $pStart + 3
$insertAt = $d.Range($insertPos, $insertPos)
$insertAt.InsertBefore(" ")

# From here on, absolute offsets after the "{m:" (3 chars) are shifted by +1
# versus the original text because of the inserted space.

# --- 2) Split "{m: " into 4 separate runs: "{" | "m" | ":" | " " ---
$s1 = $pStart + 1
$e1 = $pStart + 4
Split-Range $s1 $e1              # isolates "m: " from "{"

$s2 = $pStart + 2
$e2 = $pStart + 4
Split-Range $s2 $e2              # isolates ": " from "m"

$s3 = $pStart + 3
$e3 = $pStart + 4
Split-Range $s3 $e3              # isolates " " from ":"

# --- 3) Split "commentblock" into "comment" | "block" ---
# "comment" is 7 chars, starts right after the bookmark (pStart + 4)
$s4 = $pStart + 11
$e4 = $pStart + 16
Split-Range $s4 $e4              # isolates "block" from "comment"

# --- 4) Split " some important comment}" into " some important comment" | "}" ---
# That text is 24 chars long, starting at pStart + 16; "}" is the last char.
$s5 = $pStart + 39
$e5 = $pStart + 40
Split-Range $s5 $e5              # isolates "}" from " some important comment"

Write-Host "Final paragraph text:" $target.Range.Text
